# "starting ship fuel using item rework"
# Adds a new "Second Fuel Tank" module row to the ship-modules table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives in an Excel Table (ListObject) — grow it by one row so the
# table range/autofilter/dimension all extend together, then fill the values.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$ws.Range("B6").Value = "Second Fuel Tank "
$ws.Range("C6").Value = "Structural"
$ws.Range("D6").Value = "Adds more fuel tank size"

# Column B got a bit wider to fit the new, longer item name.
$ws.Columns.Item(2).ColumnWidth = 23.5

# Leave the selection where the editor last clicked.
$ws.Range("D30").Select()
